$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$oldText = "No one. At that moment, I just thought that while growing up everything will be great. Everyone"
$newText = "No one. At some moment in my past, I thought that while growing up everything will be great. All that hope has been shattered over time. Everyone"

$find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, $newText, 2)
